$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format first so numeric-looking strings
# (e.g. "26.953.72", "0.9995", "6.647") are stored as text, matching the
# original inlineStr/text cell type, not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.953.72"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.806.60"
$ws.Range("E3").Value = "  -1.05%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "310.36"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "0.4405"
$ws.Range("E7").Value = "  +4.19%  "
$ws.Range("D8").Value = "0.3711"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "0.07429"
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").Value = "0.8615"
$ws.Range("E10").Value = "  +2.06%  "
$ws.Range("D11").Value = "20.68"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "1.803.52"
$ws.Range("E12").Value = "  -0.53%  "
$ws.Range("D13").Value = "6.647"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "93.16"
$ws.Range("E14").Value = "  +3.70%  "
$ws.Range("D15").Value = "0.07069"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").Value = "5.282"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("E17").Value = "  -0.11%  "
$ws.Range("D18").Value = "0.000008693"
$ws.Range("E18").Value = "  -0.65%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "14.85"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("D21").Value = "26.982.77"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").Value = "5.176"
$ws.Range("E22").Value = "  +0.53%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").Value = "2.019.55"
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "1.982"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "151.13"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  -1.86%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").Value = "5.206"
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").Value = "117.53"
$ws.Range("E30").Value = "  +0.51%  "
$ws.Range("D31").Value = "0.08780"
$ws.Range("E31").Value = "  +0.47%  "
$ws.Range("D32").Value = "0.7477"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("D34").Value = "4.489"
$ws.Range("E34").Value = "  +1.46%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "0.9991"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "1.096"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "0.01972"
$ws.Range("E38").Value = "  +1.36%  "
$ws.Range("D39").Value = "0.05218"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("D41").Value = "7.080"
$ws.Range("E41").Value = "  -3.27%  "
$ws.Range("D42").Value = "2.820"
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("D43").Value = "0.1689"
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "2.144"
$ws.Range("E44").Value = "  +10.45%  "
$ws.Range("D45").Value = "8.525"
$ws.Range("D46").Value = "0.4982"
$ws.Range("E46").Value = "  +5.43%  "
$ws.Range("D47").Value = "10.42"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").Value = "104.30"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "1.672"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("B50").Value = "PaxDollar"
$ws.Range("C50").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D50").Value = "0.9989"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "0.06346"
$ws.Range("E51").Value = "  +0.32%  "

# Remove the temporary text-format styling so cells keep their original
# (default / unstyled) appearance, now that the text values are locked in.
$ws.Range("D2:D51").ClearFormats()

